$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (shifts old C,D,E -> D,E,F)
$bWidth = $ws.Range("B1").ColumnWidth
$ws.Range("C1").EntireColumn.Insert()
$ws.Range("C1").ColumnWidth = $bWidth

# New column C header duplicates the "Nam hoc" (Year) header
$ws.Range("C1").Value = "Năm học"

# Fill the new "Năm học" column with the year 2020 for both data rows
$ws.Range("C2").Value = 2020
$ws.Range("C3").Value = 2020

# Column A (So hieu / Number) updates
# Row 2: becomes a text value "49.330"
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "49.330"

# Row 3: stays numeric, value changes to 49.331
$ws.Range("A3").Value = 49.331

# Row 3 candidate changed: new name + new class code
$ws.Range("D3").Value = "Nguyen Van Nghia "
$ws.Range("E3").Value = "B12D48"

# Row 3 "Chuyen khoa" (specialty) value updates from 3 to 4
$ws.Range("F3").Value = 4

# Update selection to reflect the new last-edited cell
$ws.Range("F3").Select()
